# Auto-generated edit script: apply 2022-07-22 data updates to column I (2022 totals)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 3852
$ws.Cells.Item(3, 9).Value = 3977
$ws.Cells.Item(4, 9).Value = 929
$ws.Cells.Item(5, 9).Value = 369
$ws.Cells.Item(6, 9).Value = 4467
$ws.Cells.Item(7, 9).Value = 13594

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(6, 9).Value = 92
$ws.Cells.Item(7, 9).Value = 427
$ws.Cells.Item(8, 9).Value = 825
$ws.Cells.Item(9, 9).Value = 60
$ws.Cells.Item(10, 9).Value = 92
$ws.Cells.Item(11, 9).Value = 209
$ws.Cells.Item(14, 9).Value = 71
$ws.Cells.Item(17, 9).Value = 17
$ws.Cells.Item(19, 9).Value = 373
$ws.Cells.Item(20, 9).Value = 333
$ws.Cells.Item(23, 9).Value = 133
$ws.Cells.Item(27, 9).Value = 128
$ws.Cells.Item(29, 9).Value = 880
$ws.Cells.Item(33, 9).Value = 615
$ws.Cells.Item(36, 9).Value = 191
$ws.Cells.Item(37, 9).Value = 440
$ws.Cells.Item(42, 9).Value = 474
$ws.Cells.Item(47, 9).Value = 94
$ws.Cells.Item(48, 9).Value = 181
$ws.Cells.Item(49, 9).Value = 115
$ws.Cells.Item(50, 9).Value = 56
$ws.Cells.Item(52, 9).Value = 294
$ws.Cells.Item(53, 9).Value = 147
$ws.Cells.Item(54, 9).Value = 312
$ws.Cells.Item(55, 9).Value = 149
$ws.Cells.Item(60, 9).Value = 67
$ws.Cells.Item(61, 9).Value = 16
$ws.Cells.Item(63, 9).Value = 50
$ws.Cells.Item(64, 9).Value = 121
$ws.Cells.Item(66, 9).Value = 37
$ws.Cells.Item(67, 9).Value = 526
$ws.Cells.Item(70, 9).Value = 26
$ws.Cells.Item(71, 9).Value = 40
$ws.Cells.Item(72, 9).Value = 51
$ws.Cells.Item(73, 9).Value = 112
$ws.Cells.Item(76, 9).Value = 205
$ws.Cells.Item(77, 9).Value = 73
$ws.Cells.Item(79, 9).Value = 366
$ws.Cells.Item(83, 9).Value = 275
$ws.Cells.Item(85, 9).Value = 614
$ws.Cells.Item(89, 9).Value = 150
$ws.Cells.Item(93, 9).Value = 77
$ws.Cells.Item(94, 9).Value = 129
$ws.Cells.Item(95, 9).Value = 217
$ws.Cells.Item(99, 9).Value = 252
$ws.Cells.Item(101, 9).Value = 13594

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(6, 9).Value = 56
$ws.Cells.Item(7, 9).Value = 150

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(3, 9).Value = 22
$ws.Cells.Item(7, 9).Value = 71

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(6, 9).Value = 119
$ws.Cells.Item(7, 9).Value = 440

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 9).Value = 91
$ws.Cells.Item(7, 9).Value = 252

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 9).Value = 189
$ws.Cells.Item(5, 9).Value = 13
$ws.Cells.Item(6, 9).Value = 173
$ws.Cells.Item(7, 9).Value = 526

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 96
$ws.Cells.Item(4, 9).Value = 12
$ws.Cells.Item(6, 9).Value = 50
$ws.Cells.Item(7, 9).Value = 275

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 9).Value = 79
$ws.Cells.Item(3, 9).Value = 81
$ws.Cells.Item(7, 9).Value = 217

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 9).Value = 225
$ws.Cells.Item(5, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 615

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(4, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 115

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 9).Value = 71
$ws.Cells.Item(3, 9).Value = 62
$ws.Cells.Item(6, 9).Value = 156
$ws.Cells.Item(7, 9).Value = 312

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 262
$ws.Cells.Item(3, 9).Value = 299
$ws.Cells.Item(6, 9).Value = 242
$ws.Cells.Item(7, 9).Value = 880

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 9).Value = 139
$ws.Cells.Item(3, 9).Value = 106
$ws.Cells.Item(4, 9).Value = 17
$ws.Cells.Item(6, 9).Value = 102
$ws.Cells.Item(7, 9).Value = 373

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(4, 9).Value = 18
$ws.Cells.Item(6, 9).Value = 106
$ws.Cells.Item(7, 9).Value = 181

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(4, 9).Value = 26
$ws.Cells.Item(7, 9).Value = 205

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 159
$ws.Cells.Item(5, 9).Value = 21
$ws.Cells.Item(6, 9).Value = 151
$ws.Cells.Item(7, 9).Value = 614

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(6, 9).Value = 19
$ws.Cells.Item(7, 9).Value = 92

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 9).Value = 128
$ws.Cells.Item(7, 9).Value = 474

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 9).Value = 31
$ws.Cells.Item(7, 9).Value = 92

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 9).Value = 47
$ws.Cells.Item(3, 9).Value = 40
$ws.Cells.Item(7, 9).Value = 149

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 9).Value = 36
$ws.Cells.Item(7, 9).Value = 133

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(4, 9).Value = 23
$ws.Cells.Item(6, 9).Value = 110
$ws.Cells.Item(7, 9).Value = 366

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 9).Value = 36
$ws.Cells.Item(6, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 121

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 9).Value = 93
$ws.Cells.Item(6, 9).Value = 105
$ws.Cells.Item(7, 9).Value = 333

$ws = $wb.Worksheets.Item('Burnside')
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(7, 9).Value = 17

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(4, 9).Value = 9
$ws.Cells.Item(7, 9).Value = 191

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 9).Value = 24
$ws.Cells.Item(3, 9).Value = 21
$ws.Cells.Item(7, 9).Value = 77

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 82
$ws.Cells.Item(3, 9).Value = 100
$ws.Cells.Item(4, 9).Value = 28
$ws.Cells.Item(7, 9).Value = 294

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 9).Value = 24
$ws.Cells.Item(3, 9).Value = 23
$ws.Cells.Item(7, 9).Value = 129

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 94

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(4, 9).Value = 12
$ws.Cells.Item(7, 9).Value = 56

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 9).Value = 17
$ws.Cells.Item(7, 9).Value = 37

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 9).Value = 94
$ws.Cells.Item(6, 9).Value = 50
$ws.Cells.Item(7, 9).Value = 209

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(2, 9).Value = 23
$ws.Cells.Item(7, 9).Value = 60

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 9).Value = 40
$ws.Cells.Item(7, 9).Value = 112

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(3, 9).Value = 8
$ws.Cells.Item(7, 9).Value = 26

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 9).Value = 256
$ws.Cells.Item(3, 9).Value = 230
$ws.Cells.Item(6, 9).Value = 268
$ws.Cells.Item(7, 9).Value = 825

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 128

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(5, 9).Value = 5
$ws.Cells.Item(7, 9).Value = 67

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 9).Value = 67
$ws.Cells.Item(7, 9).Value = 147

$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(2, 9).Value = 13
$ws.Cells.Item(7, 9).Value = 40

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(6, 9).Value = 28
$ws.Cells.Item(7, 9).Value = 51

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(2, 9).Value = 20
$ws.Cells.Item(7, 9).Value = 73

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 9).Value = 148
$ws.Cells.Item(7, 9).Value = 427

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Cells.Item(6, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 16
